# Generate Report for Handoff
# Refresh the localization-status report for the new handoff run:
#  - new source GUID (68d7ba56-81d9-4787-9e59-60b168ce4a86) replaces the old one
#  - new content hash (a4378a6f70b705d7ffa793a3dceb88c08e387bad) for the generated .xlf files
#  - refreshed handoff timestamps
#  - target/handback bookkeeping cleared out since this is a brand new handoff

$wb = $excel.ActiveWorkbook

$oldGuid = "79134831-0534-4f0a-988d-4df5b37a1c1c"
$newGuid = "68d7ba56-81d9-4787-9e59-60b168ce4a86"
$newHash = "a4378a6f70b705d7ffa793a3dceb88c08e387bad"

$newHoDate = "2016-08-29 00:58:32"
$zeroDate  = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newHoDate

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$B`$2") {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-29 00:58:27"
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $zeroDate

$zhHyperlinksToDelete = @()
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$I`$2") {
        $zhHyperlinksToDelete += $hl
    }
}
foreach ($hl in $zhHyperlinksToDelete) {
    $hl.Delete()
}

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $zeroDate

$deHyperlinksToDelete = @()
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$I`$2") {
        $deHyperlinksToDelete += $hl
    }
}
foreach ($hl in $deHyperlinksToDelete) {
    $hl.Delete()
}

# ---------------------------------------------------------------------------
# Column widths for I/J narrow down on the language sheets now that the
# Target File / Handback File columns are empty
# ---------------------------------------------------------------------------
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426
$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426

Write-Host "Report refreshed for new handoff $newGuid"
